$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for column B (row => new text). Some look numeric, so we
# force them to stay text (shared-string) cells the same way Excel does
# when a cell has been pre-formatted as Text: set NumberFormat to "@"
# before typing the value, then restore the cell style to "Normal" so the
# cell itself ends up on the default style again (matching the workbook's
# original formatting, which never sets any explicit cell style).
$values = [ordered]@{
    "B2"  = "КОЙЧА"
    "B5"  = "ьЕЛЙКЙИ"
    "B6"  = "ЗъГЛАИТЪ"
    "B7"  = "ПЭТИсоТ"
    "B8"  = "1215"
    "B11" = "14"
    "B14" = "НАоЛГйТЛжЕ"
    "B16" = "34"
    "B17" = "1257"
    "B18" = "13"
    "B20" = "237"
    "B21" = "26"
    "B24" = "ДОМАГАЁТ"
    "B25" = "5"
    "B26" = "77"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Style = "Normal"
}
